# v1.0.14 updates - "Stock Type" worksheet edits
# - Add a new "Black Carbon" stock type row (inserted before "DOM: Aboveground Fast")
# - Add a new "Peat" stock type row at the end
# - Document the usage ("baseline" vs "transition triggered") of several flows
#   in the Description column for the Atmosphere rows and the Forestry Sector row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Type")

# Insert a new row for the "Black Carbon" stock type right before
# "DOM: Aboveground Fast" (original row 11), pushing the DOM rows (and
# everything below) down by one.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Black Carbon"

# Document which flows each Atmosphere pool participates in.
$ws.Range("B2").Value = "Used for baseline flows"
$ws.Range("B3").Value = "Used only for transition triggered flows"
$ws.Range("B4").Value = "Used only for transition triggered flows"
$ws.Range("B5").Value = "Used only for transition triggered flows"

# Forestry Sector row (shifted from row 21 to row 22 by the insert above)
# also gets a description.
$ws.Range("B22").Value = "Used only for transition triggered flows"

# New "Peat" stock type row appended at the end of the table.
$ws.Range("A23").Value = "Peat"

# Match the saved selection state.
$ws.Range("B30").Select()
